# Fix the UF labels for the tied "last place" rows (23-26) across the
# ranking sheets so each sheet's author-fixed database year lines up with
# the correct state abbreviation (RO/SE/AC/TO were shuffled).

$wb = $excel.ActiveWorkbook

# sheet "qtd"
$ws1 = $wb.Worksheets.Item("qtd")
$ws1.Range("A23").Value = "MT"
$ws1.Range("A24").Value = "SE"
$ws1.Range("A25").Value = "AC"
$ws1.Range("A26").Value = "RO"

# sheet "tot-arrecad"
$ws2 = $wb.Worksheets.Item("tot-arrecad")
$ws2.Range("A23").Value = "TO"
$ws2.Range("A24").Value = "SE"
$ws2.Range("A25").Value = "RO"
$ws2.Range("A26").Value = "AC"

# sheet "avg-arrecad"
$ws3 = $wb.Worksheets.Item("avg-arrecad")
$ws3.Range("A23").Value = "TO"
$ws3.Range("A24").Value = "SE"
$ws3.Range("A25").Value = "AC"
$ws3.Range("A26").Value = "RO"

# sheet "max-arrecad"
$ws4 = $wb.Worksheets.Item("max-arrecad")
$ws4.Range("A23").Value = "TO"
$ws4.Range("A24").Value = "SE"
$ws4.Range("A25").Value = "AC"
$ws4.Range("A26").Value = "RO"

# sheet "tx-sucesso"
$ws5 = $wb.Worksheets.Item("tx-sucesso")
$ws5.Range("A23").Value = "TO"
$ws5.Range("A24").Value = "SE"
$ws5.Range("A25").Value = "RO"
$ws5.Range("A26").Value = "AC"
